# Insert two new columns at C (old column C "PU 140" and D "PU 200"
# shift right to become E and F), then populate the freed-up column C
# with a "Data (bytes)" header and per-row formulas that convert the
# (now-in-column-E) PU140 value into a byte count.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("C:D").Insert() | Out-Null

$ws.Range("C1").Value = "Data (bytes)"

# C2 is entered as its own formula; C3:C21 are filled afterwards so they
# form a shared-formula group anchored at C3 (matches how Excel records
# "type formula once, then fill down over the remaining rows").
$ws.Range("C2").Formula = "=INT(1000000*E2)"
$ws.Range("C3:C21").Formula = "=INT(1000000*E3)"

$ws.Range("C2").Select() | Out-Null
